$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 467, shifting existing rows 467:508 down to 469:510.
$ws.Rows("467:468").Insert()

# Row 467: new "Primera" quality record for Apio, dated 2022-10-24 (serial 44858)
$ws.Cells.Item(467, 1).Value = 8
$ws.Cells.Item(467, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(467, 3).Value = "Coquimbo"
$ws.Cells.Item(467, 4).Value = 44858
$ws.Cells.Item(467, 5).Value = 4
$ws.Cells.Item(467, 6).Value = 100112017
$ws.Cells.Item(467, 7).Value = "Apio"
$ws.Cells.Item(467, 8).Value = "Americana (o)"
$ws.Cells.Item(467, 9).Value = "Primera"
$ws.Cells.Item(467, 10).Value = 1800
$ws.Cells.Item(467, 11).Value = 9000
$ws.Cells.Item(467, 12).Value = 10000
$ws.Cells.Item(467, 13).Value = 9500
$ws.Cells.Item(467, 14).Value = "$/docena de matas"
$ws.Cells.Item(467, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(467, 16).Value = 1583
$ws.Cells.Item(467, 17).Value = 6
$ws.Cells.Item(467, 18).Value = "Hortaliza"

# Row 468: new "Segunda" quality record for Apio, dated 2022-10-24 (serial 44858)
$ws.Cells.Item(468, 1).Value = 8
$ws.Cells.Item(468, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(468, 3).Value = "Coquimbo"
$ws.Cells.Item(468, 4).Value = 44858
$ws.Cells.Item(468, 5).Value = 4
$ws.Cells.Item(468, 6).Value = 100112017
$ws.Cells.Item(468, 7).Value = "Apio"
$ws.Cells.Item(468, 8).Value = "Americana (o)"
$ws.Cells.Item(468, 9).Value = "Segunda"
$ws.Cells.Item(468, 10).Value = 1200
$ws.Cells.Item(468, 11).Value = 7000
$ws.Cells.Item(468, 12).Value = 8000
$ws.Cells.Item(468, 13).Value = 7500
$ws.Cells.Item(468, 14).Value = "$/docena de matas"
$ws.Cells.Item(468, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(468, 16).Value = 1250
$ws.Cells.Item(468, 17).Value = 6
$ws.Cells.Item(468, 18).Value = "Hortaliza"
